$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "30.334.04"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -2.88%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.935.69"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -3.02%  "
$ws.Range("E4").Value = "  +0.07%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "251.01"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.7246"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -6.67%  "
$ws.Range("E7").Value = "  +0.12%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3311"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -4.80%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "28.06"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.73%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07195"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.82%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.8114"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -3.39%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.08093"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.33%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.937.11"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -2.93%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.494"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.32%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "94.65"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -6.09%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "15.26"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.47%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "30.348.92"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -2.80%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000008306"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +3.81%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "250.58"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -7.93%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "5.917"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "2.192.14"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.75%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  +0.08%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "6.996"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.32%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "9.755"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -2.35%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "163.69"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.59%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.378"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.84%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "19.29"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -2.93%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.1327"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -7.00%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.569"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.52%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.348"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.09%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.445"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -3.50%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.181"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -5.88%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.05203"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("E35").Value = "  +3.28%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7514"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.746"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.82%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01984"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.90%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.836"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.66%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "79.85"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -4.04%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "6.441"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -4.90%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.4541"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("E43").Value = "  -4.65%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.8490"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.66%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.06%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "101.99"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -2.86%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "9.803"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.35%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.490"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -3.12%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "36.83"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -1.97%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.4185"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -3.33%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.06040"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
